$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Michael's row: NRIC (B2) and Password (E2) are both overwritten with "c"
# while testing the register/login flow for an Officer account.
$ws.Range("B2").Value = "c"
$ws.Range("E2").Value = "c"

# Move the active selection to F2, matching the saved selection state.
$ws.Range("F2").Select()
